$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "1.00" or "46.30");
# force text format so Excel does not silently coerce/normalize these values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.695.53"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.487.21"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.96"
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.08"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.484.18"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.57"
$ws.Range("E11").Value = "  +5.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.424"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000214"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.075.74"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.54"
$ws.Range("E15").Value = "  -3.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.485.56"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.801.74"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("E19").Value = "  -3.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.38"
$ws.Range("E20").Value = "  -1.81%  "
$ws.Range("E21").Value = "  +3.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "439.18"
$ws.Range("E22").Value = "  -2.18%  "
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.62"
$ws.Range("E24").Value = "  +2.59%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.622.83"
$ws.Range("E26").Value = "  -1.52%  "
$ws.Range("E27").Value = "  -6.72%  "
$ws.Range("E28").Value = "  -5.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.39"
$ws.Range("E29").Value = "  -5.62%  "
$ws.Range("E30").Value = "  -1.77%  "
$ws.Range("E31").Value = "  -3.62%  "
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.42"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.05"
$ws.Range("E35").Value = "  -4.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.476.68"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").Value = "  -4.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.93"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "177.49"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0889"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("E43").Value = "  -8.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.42"
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.11"
$ws.Range("E46").Value = "  -3.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.30"
$ws.Range("E47").Value = "  +1.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.23"
$ws.Range("E48").Value = "  -6.43%  "
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("E50").Value = "  -7.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.983"
$ws.Range("E51").Value = "  -2.44%  "
